$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvestorKyc")

# Add the new "Agreement Committed Amount" column in R
$ws.Range("R1").Value = "Agreement Committed Amount"
$ws.Range("R2").Value = 1000000
$ws.Range("R3").Value = 2000000
$ws.Range("R4").Value = 3000000
$ws.Range("R5").Value = 4000000

# Rename "Full Name" header (B1) to "Investing Entity"
$ws.Range("B1").Value = "Investing Entity"

# Update the active selection to match the saved view state
$ws.Range("R6").Select()
